$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename shared string used by column AD header: "apya" -> "zahid ur Rehman" ---
# (this also appends "zahid ur Rehman" as a new shared string once "apya" becomes unused)
$ws.Range("AD1").Value = "zahid ur Rehman"

# --- Update "shoaib" (column L) values for a few days ---
$ws.Range("L2").Value = 2
$ws.Range("L25").Value = 2
$ws.Range("L26").Value = 2

$arr27 = New-Object "object[,]" 1,43
$arr27[0,0] = 0.5
$arr27[0,1] = 3
$arr27[0,2] = 1.5
$arr27[0,3] = 2
$arr27[0,4] = 2
$arr27[0,5] = 2
$arr27[0,6] = 3.5
$arr27[0,7] = 2.5
$arr27[0,8] = 1
$arr27[0,9] = 5
$arr27[0,10] = 2
$arr27[0,11] = 3
$arr27[0,12] = 3
$arr27[0,13] = 0
$arr27[0,14] = 1.5
$arr27[0,15] = 1
$arr27[0,16] = 1
$arr27[0,17] = 1
$arr27[0,18] = 1
$arr27[0,19] = 1.5
$arr27[0,20] = 5
$arr27[0,21] = 2.5
$arr27[0,22] = 0
$arr27[0,23] = 0.5
$arr27[0,24] = 2
$arr27[0,25] = 2
$arr27[0,26] = 1
$arr27[0,27] = 0
$arr27[0,28] = 1.5
$arr27[0,29] = 2
$arr27[0,30] = 3
$arr27[0,31] = 3.5
$arr27[0,32] = 1.5
$arr27[0,33] = 0
$arr27[0,34] = 1
$arr27[0,35] = 1.5
$arr27[0,36] = 1.5
$arr27[0,37] = 2
$arr27[0,38] = 0
$arr27[0,39] = 2.5
$arr27[0,40] = 0.5
$arr27[0,41] = 2.5
$arr27[0,42] = 0
$ws.Range("H27:AX27").Value = $arr27

$arr28 = New-Object "object[,]" 1,43
$arr28[0,0] = 0.5
$arr28[0,1] = 3
$arr28[0,2] = 1.5
$arr28[0,3] = 2
$arr28[0,4] = 2
$arr28[0,5] = 2
$arr28[0,6] = 3.5
$arr28[0,7] = 2.5
$arr28[0,8] = 1
$arr28[0,9] = 5
$arr28[0,10] = 2
$arr28[0,11] = 3
$arr28[0,12] = 3
$arr28[0,13] = 0
$arr28[0,14] = 1.5
$arr28[0,15] = 1
$arr28[0,16] = 1
$arr28[0,17] = 1
$arr28[0,18] = 1
$arr28[0,19] = 1.5
$arr28[0,20] = 5
$arr28[0,21] = 2.5
$arr28[0,22] = 0
$arr28[0,23] = 0.5
$arr28[0,24] = 2
$arr28[0,25] = 2
$arr28[0,26] = 1
$arr28[0,27] = 0
$arr28[0,28] = 1.5
$arr28[0,29] = 2
$arr28[0,30] = 3
$arr28[0,31] = 3.5
$arr28[0,32] = 1.5
$arr28[0,33] = 0
$arr28[0,34] = 1
$arr28[0,35] = 1.5
$arr28[0,36] = 1.5
$arr28[0,37] = 0
$arr28[0,38] = 0
$arr28[0,39] = 2.5
$arr28[0,40] = 0.5
$arr28[0,41] = 2.5
$arr28[0,42] = 0
$ws.Range("H28:AX28").Value = $arr28

$arr29 = New-Object "object[,]" 1,43
$arr29[0,0] = 0.5
$arr29[0,1] = 3
$arr29[0,2] = 1.5
$arr29[0,3] = 2
$arr29[0,4] = 2
$arr29[0,5] = 2
$arr29[0,6] = 3.5
$arr29[0,7] = 2.5
$arr29[0,8] = 1
$arr29[0,9] = 5
$arr29[0,10] = 2
$arr29[0,11] = 3
$arr29[0,12] = 3
$arr29[0,13] = 0
$arr29[0,14] = 1.5
$arr29[0,15] = 1
$arr29[0,16] = 1
$arr29[0,17] = 1
$arr29[0,18] = 1
$arr29[0,19] = 1.5
$arr29[0,20] = 5
$arr29[0,21] = 2.5
$arr29[0,22] = 0
$arr29[0,23] = 0.5
$arr29[0,24] = 2
$arr29[0,25] = 2
$arr29[0,26] = 1
$arr29[0,27] = 0
$arr29[0,28] = 1.5
$arr29[0,29] = 2
$arr29[0,30] = 0
$arr29[0,31] = 3.5
$arr29[0,32] = 1.5
$arr29[0,33] = 0
$arr29[0,34] = 1
$arr29[0,35] = 1.5
$arr29[0,36] = 1.5
$arr29[0,37] = 0
$arr29[0,38] = 0
$arr29[0,39] = 2.5
$arr29[0,40] = 0.5
$arr29[0,41] = 2.5
$arr29[0,42] = 0
$ws.Range("H29:AX29").Value = $arr29

$arr30 = New-Object "object[,]" 1,43
$arr30[0,0] = 0.5
$arr30[0,1] = 3
$arr30[0,2] = 1.5
$arr30[0,3] = 2
$arr30[0,4] = 2
$arr30[0,5] = 2
$arr30[0,6] = 3.5
$arr30[0,7] = 2.5
$arr30[0,8] = 1
$arr30[0,9] = 5
$arr30[0,10] = 2
$arr30[0,11] = 3
$arr30[0,12] = 3
$arr30[0,13] = 0
$arr30[0,14] = 1.5
$arr30[0,15] = 1
$arr30[0,16] = 1
$arr30[0,17] = 1
$arr30[0,18] = 1
$arr30[0,19] = 1.5
$arr30[0,20] = 5
$arr30[0,21] = 2.5
$arr30[0,22] = 0
$arr30[0,23] = 0.5
$arr30[0,24] = 2
$arr30[0,25] = 2
$arr30[0,26] = 1
$arr30[0,27] = 0
$arr30[0,28] = 1.5
$arr30[0,29] = 2
$arr30[0,30] = 0
$arr30[0,31] = 3.5
$arr30[0,32] = 1.5
$arr30[0,33] = 0
$arr30[0,34] = 1
$arr30[0,35] = 1.5
$arr30[0,36] = 1.5
$arr30[0,37] = 0
$arr30[0,38] = 0
$arr30[0,39] = 2.5
$arr30[0,40] = 0.5
$arr30[0,41] = 2.5
$arr30[0,42] = 0
$ws.Range("H30:AX30").Value = $arr30

$arr31 = New-Object "object[,]" 1,43
$arr31[0,0] = 0.5
$arr31[0,1] = 3
$arr31[0,2] = 1.5
$arr31[0,3] = 2
$arr31[0,4] = 2
$arr31[0,5] = 2
$arr31[0,6] = 3.5
$arr31[0,7] = 2.5
$arr31[0,8] = 1
$arr31[0,9] = 5
$arr31[0,10] = 2
$arr31[0,11] = 3
$arr31[0,12] = 3
$arr31[0,13] = 0
$arr31[0,14] = 1.5
$arr31[0,15] = 1
$arr31[0,16] = 1
$arr31[0,17] = 1
$arr31[0,18] = 1
$arr31[0,19] = 0
$arr31[0,20] = 5
$arr31[0,21] = 2.5
$arr31[0,22] = 2
$arr31[0,23] = 0.5
$arr31[0,24] = 2
$arr31[0,25] = 2
$arr31[0,26] = 1
$arr31[0,27] = 1.5
$arr31[0,28] = 1.5
$arr31[0,29] = 2
$arr31[0,30] = 3
$arr31[0,31] = 3.5
$arr31[0,32] = 1.5
$arr31[0,33] = 2
$arr31[0,34] = 1
$arr31[0,35] = 1.5
$arr31[0,36] = 1.5
$arr31[0,37] = 2
$arr31[0,38] = 0
$arr31[0,39] = 2.5
$arr31[0,40] = 0.5
$arr31[0,41] = 2.5
$arr31[0,42] = 0
$ws.Range("H31:AX31").Value = $arr31


# --- AH29:AH31 previously had a special "no value" fill style; now that they carry
#     real data, clear the leftover yellow fill so they match the other data cells ---
$ws.Range("AH29:AH31").Interior.Pattern = -4142

# --- Update the saved view: scroll back to A1 and select H27 ---
$ws.Range("H27").Select()
